# Applies the "revised ppt template formats" edit:
#  - merges the split Chinese runs "用户"+"界面" -> "用户界面" (User Interface box)
#  - merges the split Chinese runs "用户"+"管理" -> "用户管理" (Account Administration box)
#  - strips the text drop-shadow effect that was baked into every label on the slide

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cr = [char]13

# --- Title shape ("Components") ---------------------------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Font.Shadow = $false

# --- The big group of rectangle labels ---------------------------------
$grp = $s.Shapes.Item(2)

for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $shp = $grp.GroupItems.Item($i)
    if ($shp.HasTextFrame) {
        $shp.TextFrame.TextRange.Font.Shadow = $false
    }

    if ($shp.Name -eq "Rectangle 5") {
        # "User Interface" / 用户界面
        $shp.TextFrame.TextRange.Text = "User Interface" + $cr + "用户界面"
        $shp.TextFrame.TextRange.Font.Shadow = $false
    }
    elseif ($shp.Name -eq "Rectangle 7") {
        # "Account Administration" / 用户管理
        $shp.TextFrame.TextRange.Text = "Account Administration" + $cr + "用户管理"
        $shp.TextFrame.TextRange.Font.Shadow = $false
    }
}
